$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: copy formatting (thin border style, s="4") from the last existing
# data row (1243) down across the new row range so every new cell inherits
# the same style used throughout the table, without creating new style entries.
$fmtSrc = $ws.Range("A1243:N1243")
$fmtDst = $ws.Range("A1244:N1272")
$fmtSrc.Copy($fmtDst)

# Step 2: write the actual cell values for the new rows.
# Row 1244
$ws.Range("A1244").Value = 'Unit1'
$ws.Range("B1244").Formula = '="2019-04-16"'
$ws.Range("C1244").Value = 1212
$ws.Range("D1244").Value = 12
$ws.Range("E1244").Value = 12
$ws.Range("F1244").Value = 12
$ws.Range("G1244").Value = 34
$ws.Range("H1244").Value = 34
$ws.Range("I1244").Value = 34
$ws.Range("J1244").Value = 34
$ws.Range("K1244").Value = 34
$ws.Range("L1244").Value = 34
$ws.Range("M1244").Value = 34
$ws.Range("N1244").Value = 34

# Row 1245
$ws.Range("A1245").Value = 'Unit1'
$ws.Range("B1245").Formula = '="2019-04-16"'
$ws.Range("C1245").Value = 1122222
$ws.Range("D1245").Value = 145
$ws.Range("E1245").Value = 125
$ws.Range("F1245").Value = 125
$ws.Range("G1245").Value = 126
$ws.Range("H1245").Value = 127
$ws.Range("I1245").Value = 1281
$ws.Range("J1245").Value = 129
$ws.Range("K1245").Value = 130
$ws.Range("L1245").Value = 140
$ws.Range("M1245").Value = 150
$ws.Range("N1245").Value = 160

# Row 1246
$ws.Range("A1246").Value = 'Unit1'
$ws.Range("B1246").Formula = '="2019-04-01"'
$ws.Range("C1246").Value = 12
$ws.Range("D1246").Value = 12
$ws.Range("E1246").Value = 12
$ws.Range("F1246").Value = 12
$ws.Range("G1246").Value = 12
$ws.Range("H1246").Value = 12
$ws.Range("I1246").Value = 12
$ws.Range("J1246").Value = 12
$ws.Range("K1246").Value = 12
$ws.Range("L1246").Value = 12
$ws.Range("M1246").Value = 12
$ws.Range("N1246").Value = 12

# Row 1247
$ws.Range("A1247").Value = 'Unit1'
$ws.Range("B1247").Formula = '="2019-04-16"'
$ws.Range("C1247").Value = 122
$ws.Range("D1247").Value = 12
$ws.Range("E1247").Value = 112
$ws.Range("F1247").Value = 1212
$ws.Range("G1247").Value = 212
$ws.Range("H1247").Value = 1212
$ws.Range("I1247").Value = 1212
$ws.Range("J1247").Value = 1212
$ws.Range("K1247").Value = 12
$ws.Range("L1247").Value = 12
$ws.Range("M1247").Value = 12
$ws.Range("N1247").Value = 12

# Row 1248
$ws.Range("A1248").Value = 'Unit1'
$ws.Range("B1248").Formula = '="2019-04-09"'
$ws.Range("C1248").Value = 12
$ws.Range("D1248").Value = 12
$ws.Range("E1248").Value = 12
$ws.Range("F1248").Value = 12
$ws.Range("G1248").Value = 12
$ws.Range("H1248").Value = 12
$ws.Range("I1248").Value = 12
$ws.Range("J1248").Value = 12
$ws.Range("K1248").Value = 12
$ws.Range("L1248").Value = 12
$ws.Range("M1248").Value = 12
$ws.Range("N1248").Value = 12

# Row 1249
$ws.Range("A1249").Value = 'Unit1'
$ws.Range("B1249").Formula = '="2019-04-15"'
$ws.Range("C1249").Value = 12
$ws.Range("D1249").Value = 12
$ws.Range("E1249").Value = 12
$ws.Range("F1249").Value = 12
$ws.Range("G1249").Value = 12
$ws.Range("H1249").Value = 12
$ws.Range("I1249").Value = 12
$ws.Range("J1249").Value = 12
$ws.Range("K1249").Value = 12
$ws.Range("L1249").Value = 12
$ws.Range("M1249").Value = 12
$ws.Range("N1249").Value = 12

# Row 1250
$ws.Range("A1250").Value = 'Unit1'
$ws.Range("B1250").Formula = '="2019-04-16"'
$ws.Range("C1250").Value = 121
$ws.Range("D1250").Value = 12
$ws.Range("E1250").Value = 12
$ws.Range("F1250").Value = 12
$ws.Range("G1250").Value = 12
$ws.Range("H1250").Value = 12
$ws.Range("I1250").Value = 12
$ws.Range("J1250").Value = 12
$ws.Range("K1250").Value = 12
$ws.Range("L1250").Value = 12
$ws.Range("M1250").Value = 12
$ws.Range("N1250").Value = 12

# Row 1251
$ws.Range("A1251").Value = 'Unit1'
$ws.Range("B1251").Formula = '="2019-04-16"'
$ws.Range("C1251").Value = 12
$ws.Range("D1251").Value = 12
$ws.Range("E1251").Value = 12
$ws.Range("F1251").Value = 12
$ws.Range("G1251").Value = 112
$ws.Range("H1251").Value = 12
$ws.Range("I1251").Value = 12
$ws.Range("J1251").Value = 12
$ws.Range("K1251").Value = 12
$ws.Range("L1251").Value = 12
$ws.Range("M1251").Value = 12
$ws.Range("N1251").Value = 21

# Row 1252
$ws.Range("A1252").Value = 'Unit1'
$ws.Range("B1252").Formula = '="2019-04-16"'
$ws.Range("C1252").Value = 12
$ws.Range("D1252").Value = 'qq'
$ws.Range("E1252").Value = 12
$ws.Range("F1252").Value = 12
$ws.Range("G1252").Value = 112
$ws.Range("H1252").Value = 12
$ws.Range("I1252").Value = 12
$ws.Range("J1252").Value = 12
$ws.Range("K1252").Value = 12
$ws.Range("L1252").Value = 12
$ws.Range("M1252").Value = 12
$ws.Range("N1252").Value = 12

# Row 1253
$ws.Range("A1253").Value = 'Unit1'
$ws.Range("B1253").Formula = '="2019-04-16"'
$ws.Range("C1253").Value = 1212
$ws.Range("D1253").Value = 12213
$ws.Range("E1253").Value = 12312
$ws.Range("F1253").Value = 3123
$ws.Range("G1253").Value = 213
$ws.Range("H1253").Value = 21312
$ws.Range("I1253").Value = 3123123
$ws.Range("J1253").Value = 213
$ws.Range("K1253").Value = 123
$ws.Range("L1253").Value = 123
$ws.Range("M1253").Value = 123
$ws.Range("N1253").Value = 123

# Row 1254
$ws.Range("A1254").Value = 'Unit1'
$ws.Range("B1254").Formula = '="2019-04-16"'
$ws.Range("C1254").Value = 12
$ws.Range("D1254").Value = 12
$ws.Range("E1254").Value = 12
$ws.Range("F1254").Value = 12
$ws.Range("G1254").Value = 12
$ws.Range("H1254").Value = 12
$ws.Range("I1254").Value = 12
$ws.Range("J1254").Value = 12
$ws.Range("K1254").Value = 12
$ws.Range("L1254").Value = 12
$ws.Range("M1254").Value = 12
$ws.Range("N1254").Value = 12

# Row 1255
$ws.Range("A1255").Value = 'Unit1'
$ws.Range("B1255").Formula = '="2019-04-16"'
$ws.Range("C1255").Value = 12
$ws.Range("D1255").Value = 12
$ws.Range("E1255").Value = 12
$ws.Range("F1255").Value = 12
$ws.Range("G1255").Value = 12
$ws.Range("H1255").Value = 12
$ws.Range("I1255").Value = 12
$ws.Range("J1255").Value = 12
$ws.Range("K1255").Value = 12
$ws.Range("L1255").Value = 12
$ws.Range("M1255").Value = 12
$ws.Range("N1255").Value = 12

# Row 1256
$ws.Range("A1256").Value = 'Unit1'
$ws.Range("B1256").Formula = '="2019-04-16"'
$ws.Range("C1256").Value = 1121
$ws.Range("D1256").Value = 12
$ws.Range("E1256").Value = 12
$ws.Range("F1256").Value = 12
$ws.Range("G1256").Value = 12
$ws.Range("H1256").Value = 12
$ws.Range("I1256").Value = 12
$ws.Range("J1256").Value = 112
$ws.Range("K1256").Value = 121
$ws.Range("L1256").Value = 12
$ws.Range("M1256").Value = 12
$ws.Range("N1256").Value = 12

# Row 1257
$ws.Range("A1257").Value = 'Unit1'
$ws.Range("B1257").Formula = '="2019-04-16"'
$ws.Range("C1257").Value = 12
$ws.Range("D1257").Value = 12
$ws.Range("E1257").Value = 12
$ws.Range("F1257").Value = 12
$ws.Range("G1257").Value = 12
$ws.Range("H1257").Value = 12
$ws.Range("I1257").Value = 12
$ws.Range("J1257").Value = 12
$ws.Range("K1257").Value = 12
$ws.Range("L1257").Value = 12
$ws.Range("M1257").Value = 12
$ws.Range("N1257").Value = 12

# Row 1258
$ws.Range("A1258").Value = 'Unit1'
$ws.Range("B1258").Formula = '="2019-04-16"'
$ws.Range("C1258").Value = 12
$ws.Range("D1258").Value = 123
$ws.Range("E1258").Value = 123
$ws.Range("F1258").Value = 123
$ws.Range("G1258").Value = 123
$ws.Range("H1258").Value = 123
$ws.Range("I1258").Value = 123
$ws.Range("J1258").Value = 123
$ws.Range("K1258").Value = 213
$ws.Range("L1258").Value = 123
$ws.Range("M1258").Value = 123
$ws.Range("N1258").Value = 123

# Row 1259
$ws.Range("A1259").Value = 'Unit1'
$ws.Range("B1259").Formula = '="2019-04-23"'
$ws.Range("C1259").Value = 123
$ws.Range("D1259").Value = 123
$ws.Range("E1259").Value = 123
$ws.Range("F1259").Value = 1233
$ws.Range("G1259").Value = 123
$ws.Range("H1259").Value = 213
$ws.Range("I1259").Value = 123
$ws.Range("J1259").Value = 123
$ws.Range("K1259").Value = 123
$ws.Range("L1259").Value = 123
$ws.Range("M1259").Value = 123
$ws.Range("N1259").Value = 123

# Row 1260
$ws.Range("A1260").Value = 'Unit1'
$ws.Range("B1260").Formula = '="2019-04-16"'
$ws.Range("C1260").Value = 12
$ws.Range("D1260").Value = 123
$ws.Range("E1260").Value = 123
$ws.Range("F1260").Value = 123
$ws.Range("G1260").Value = 123
$ws.Range("H1260").Value = 123
$ws.Range("I1260").Value = 123
$ws.Range("J1260").Value = 123
$ws.Range("K1260").Value = 123
$ws.Range("L1260").Value = 123
$ws.Range("M1260").Value = 123
$ws.Range("N1260").Value = 123

# Row 1261
$ws.Range("A1261").Value = 'Unit1'
$ws.Range("B1261").Formula = '="2019-04-16"'
$ws.Range("C1261").Value = 213
$ws.Range("D1261").Value = 213
$ws.Range("E1261").Value = 123
$ws.Range("F1261").Value = 123
$ws.Range("G1261").Value = 123
$ws.Range("H1261").Value = 123
$ws.Range("I1261").Value = 123
$ws.Range("J1261").Value = 123
$ws.Range("K1261").Value = 123
$ws.Range("L1261").Value = 123
$ws.Range("M1261").Value = 123
$ws.Range("N1261").Value = 123

# Row 1262
$ws.Range("A1262").Value = 'Unit1'
$ws.Range("B1262").Formula = '="2019-04-16"'
$ws.Range("C1262").Value = 56
$ws.Range("D1262").Value = 'ghfgh'
$ws.Range("E1262").Value = 5656
$ws.Range("F1262").Value = 56
$ws.Range("G1262").Value = 67
$ws.Range("H1262").Value = 67
$ws.Range("I1262").Value = 67
$ws.Range("J1262").Value = 56
$ws.Range("K1262").Value = 56
$ws.Range("L1262").Value = 56
$ws.Range("M1262").Value = 65
$ws.Range("N1262").Value = '7fghj'

# Row 1263
$ws.Range("A1263").Value = 'Unit1'
$ws.Range("B1263").Formula = '="2019-04-16"'
$ws.Range("C1263").Value = 787
$ws.Range("D1263").Value = 78
$ws.Range("E1263").Value = 787
$ws.Range("F1263").Value = 787
$ws.Range("G1263").Value = 78
$ws.Range("H1263").Value = 78
$ws.Range("I1263").Value = 78
$ws.Range("J1263").Value = 78
$ws.Range("K1263").Value = 78
$ws.Range("L1263").Value = 78
$ws.Range("M1263").Value = 78
$ws.Range("N1263").Value = '786ghj'

# Row 1264
$ws.Range("A1264").Value = 'Unit1'
$ws.Range("B1264").Formula = '="2019-04-03"'
$ws.Range("C1264").Value = 565
$ws.Range("D1264").Value = 67
$ws.Range("E1264").Value = 67
$ws.Range("F1264").Value = 6767
$ws.Range("G1264").Value = 6767
$ws.Range("H1264").Value = 67
$ws.Range("I1264").Value = 67
$ws.Range("J1264").Value = 67
$ws.Range("K1264").Value = 89
$ws.Range("L1264").Value = 899
$ws.Range("M1264").Value = 89
$ws.Range("N1264").Value = 'hi'

# Row 1265
$ws.Range("A1265").Value = 'Unit1'
$ws.Range("B1265").Formula = '="2019-04-16"'
$ws.Range("C1265").Value = 778
$ws.Range("D1265").Value = 787
$ws.Range("E1265").Value = 878
$ws.Range("F1265").Value = 78
$ws.Range("G1265").Value = 78
$ws.Range("H1265").Value = 78
$ws.Range("I1265").Value = 78
$ws.Range("J1265").Value = 78
$ws.Range("K1265").Value = 78
$ws.Range("L1265").Value = 78
$ws.Range("M1265").Value = 78
$ws.Range("N1265").Value = 'hjhjk'

# Row 1266
$ws.Range("A1266").Value = 'Unit1'
$ws.Range("B1266").Formula = '="2019-04-09"'
$ws.Range("C1266").Value = 321
$ws.Range("D1266").Value = 123
$ws.Range("E1266").Value = 12312
$ws.Range("F1266").Value = 3123
$ws.Range("G1266").Value = 123
$ws.Range("H1266").Value = 123
$ws.Range("I1266").Value = 123
$ws.Range("J1266").Value = 123
$ws.Range("K1266").Value = 123
$ws.Range("L1266").Value = 123
$ws.Range("M1266").Value = 123
$ws.Range("N1266").Value = 123

# Row 1267
$ws.Range("A1267").Value = 'Unit1'
$ws.Range("B1267").Formula = '="2019-04-16"'
$ws.Range("C1267").Value = 1221
$ws.Range("D1267").Value = 12
$ws.Range("E1267").Value = 12
$ws.Range("F1267").Value = 12
$ws.Range("G1267").Value = 12
$ws.Range("H1267").Value = 12
$ws.Range("I1267").Value = 12
$ws.Range("J1267").Value = 12
$ws.Range("K1267").Value = 12
$ws.Range("L1267").Value = 12
$ws.Range("M1267").Value = 12
$ws.Range("N1267").Value = 12

# Row 1268
$ws.Range("A1268").Value = 'Unit1'
$ws.Range("B1268").Formula = '="2019-04-16"'
$ws.Range("C1268").Value = 23123
$ws.Range("D1268").Value = 123
$ws.Range("E1268").Value = 213
$ws.Range("F1268").Value = 123
$ws.Range("G1268").Value = 213
$ws.Range("H1268").Value = 123
$ws.Range("I1268").Value = 123
$ws.Range("J1268").Value = 123
$ws.Range("K1268").Value = 213
$ws.Range("L1268").Value = 123
$ws.Range("M1268").Value = 123
$ws.Range("N1268").Value = 123

# Row 1269
$ws.Range("A1269").Value = 'Unit3'
$ws.Range("B1269").Formula = '="2019-04-17"'
$ws.Range("C1269").Value = 1212
$ws.Range("D1269").Value = 1212
$ws.Range("E1269").Value = 12
$ws.Range("F1269").Value = 12
$ws.Range("G1269").Value = 12
$ws.Range("H1269").Value = 12
$ws.Range("I1269").Value = 12
$ws.Range("J1269").Value = 12
$ws.Range("K1269").Value = 12
$ws.Range("L1269").Value = 12
$ws.Range("M1269").Value = 12
$ws.Range("N1269").Value = 12

# Row 1270
$ws.Range("A1270").Value = 'Unit3'
$ws.Range("B1270").Formula = '="2019-04-09"'
$ws.Range("C1270").Value = 123123
$ws.Range("D1270").Value = 123
$ws.Range("E1270").Value = 123
$ws.Range("F1270").Value = 123
$ws.Range("G1270").Value = 123
$ws.Range("H1270").Value = 123
$ws.Range("I1270").Value = 213
$ws.Range("J1270").Value = 123
$ws.Range("K1270").Value = 123
$ws.Range("L1270").Value = 123
$ws.Range("M1270").Value = 123
$ws.Range("N1270").Value = 123

# Row 1271
$ws.Range("A1271").Value = 'Unit3'
$ws.Range("B1271").Formula = '="2019-04-10"'
$ws.Range("C1271").Value = 21312321
$ws.Range("D1271").Value = 321312
$ws.Range("E1271").Value = 3123
$ws.Range("F1271").Value = 123123
$ws.Range("G1271").Value = 21312
$ws.Range("H1271").Value = 3123
$ws.Range("I1271").Value = 123
$ws.Range("J1271").Value = 123
$ws.Range("K1271").Value = 123
$ws.Range("L1271").Value = 123
$ws.Range("M1271").Value = 123
$ws.Range("N1271").Value = 12312

# Row 1272
$ws.Range("A1272").Value = 'Unit3'
$ws.Range("B1272").Formula = '="2019-04-18"'
$ws.Range("C1272").Value = 123
$ws.Range("D1272").Value = 123
$ws.Range("E1272").Value = 123123
$ws.Range("F1272").Value = 213
$ws.Range("G1272").Value = 12312
$ws.Range("H1272").Value = 3123
$ws.Range("I1272").Value = 123
$ws.Range("J1272").Value = 123
$ws.Range("K1272").Value = 123
$ws.Range("L1272").Value = 123
$ws.Range("M1272").Value = 123
$ws.Range("N1272").Value = 123

# Step 3: the Formula assignments above produced live formulas for the
# date-like text cells (so Excel would not silently reinterpret them as
# real dates/serial numbers). Convert those formula cells to plain literal
# text values in place, preserving the existing border style.
$dateCells = @(
    "B1244",
    "B1245",
    "B1246",
    "B1247",
    "B1248",
    "B1249",
    "B1250",
    "B1251",
    "B1252",
    "B1253",
    "B1254",
    "B1255",
    "B1256",
    "B1257",
    "B1258",
    "B1259",
    "B1260",
    "B1261",
    "B1262",
    "B1263",
    "B1264",
    "B1265",
    "B1266",
    "B1267",
    "B1268",
    "B1269",
    "B1270",
    "B1271",
    "B1272"
)
foreach ($addr in $dateCells) {
    $cell = $ws.Range($addr)
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
$excel.CutCopyMode = 0